# Auto-generated Excel COM-interop script
# Applies the weekly data rotation described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44335
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10500
$ws.Range("S2").Value = 583

# Row 3
$ws.Range("D3").Value = 44335
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("S3").Value = 500

# Row 4
$ws.Range("D4").Value = 44189
$ws.Range("K4").Value = 'Red Beaut'
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = '$/caja 15 kilos granel'
$ws.Range("S4").Value = 833
$ws.Range("T4").Value = 15

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("K5").Value = 'Red Beaut'
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44202
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("S6").Value = 806

# Row 7
$ws.Range("D7").Value = 44202
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 667

# Row 8
$ws.Range("D8").Value = 44236
$ws.Range("K8").Value = 'Lemon'
$ws.Range("Q8").Value = '$/caja 16 kilos granel'
$ws.Range("S8").Value = 906
$ws.Range("T8").Value = 16

# Row 9
$ws.Range("D9").Value = 44236
$ws.Range("K9").Value = 'Lemon'
$ws.Range("Q9").Value = '$/caja 16 kilos granel'
$ws.Range("S9").Value = 750
$ws.Range("T9").Value = 16

# Row 10
$ws.Range("D10").Value = 44299
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12500
$ws.Range("S10").Value = 694

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 11000
$ws.Range("S11").Value = 611

# Row 12
$ws.Range("D12").Value = 44218
$ws.Range("K12").Value = 'Black Amber'
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 11000
$ws.Range("P12").Value = 10500
$ws.Range("Q12").Value = '$/caja 16 kilos granel'
$ws.Range("S12").Value = 656
$ws.Range("T12").Value = 16

# Row 13
$ws.Range("D13").Value = 44218
$ws.Range("K13").Value = 'Black Amber'
$ws.Range("N13").Value = 9000
$ws.Range("O13").Value = 9000
$ws.Range("P13").Value = 9000
$ws.Range("Q13").Value = '$/caja 16 kilos granel'
$ws.Range("S13").Value = 562
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 44223
$ws.Range("K14").Value = 'Black Amber'
$ws.Range("O14").Value = 11000
$ws.Range("P14").Value = 10500
$ws.Range("S14").Value = 656

# Row 15
$ws.Range("D15").Value = 44223
$ws.Range("K15").Value = 'Black Amber'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 9000
$ws.Range("P15").Value = 9000
$ws.Range("S15").Value = 562

# Row 16
$ws.Range("D16").Value = 44246
$ws.Range("K16").Value = 'Angeleno'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("S16").Value = 625

# Row 17
$ws.Range("D17").Value = 44246
$ws.Range("K17").Value = 'Angeleno'
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 500

# Row 18
$ws.Range("D18").Value = 44307
$ws.Range("M18").Value = 200

# Row 19
$ws.Range("D19").Value = 44307
$ws.Range("M19").Value = 100

# Row 20
$ws.Range("D20").Value = 44328
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 9500
$ws.Range("S20").Value = 528

# Row 21
$ws.Range("D21").Value = 44328
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 444

# Row 22
$ws.Range("D22").Value = 44285
$ws.Range("N22").Value = 9000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 9500
$ws.Range("Q22").Value = '$/caja 18 kilos granel'
$ws.Range("S22").Value = 528

# Row 23
$ws.Range("D23").Value = 44285
$ws.Range("N23").Value = 8000
$ws.Range("O23").Value = 8000
$ws.Range("P23").Value = 8000
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("S23").Value = 444

# Row 24
$ws.Range("D24").Value = 44266
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 9500
$ws.Range("Q24").Value = '$/caja 18 kilos granel'
$ws.Range("S24").Value = 528
$ws.Range("T24").Value = 18

# Row 25
$ws.Range("D25").Value = 44266
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 8000
$ws.Range("P25").Value = 8000
$ws.Range("Q25").Value = '$/caja 18 kilos granel'
$ws.Range("S25").Value = 444
$ws.Range("T25").Value = 18

# Row 26
$ws.Range("D26").Value = 44343
$ws.Range("M26").Value = 200

# Row 27
$ws.Range("D27").Value = 44343
$ws.Range("M27").Value = 100

# Row 28
$ws.Range("D28").Value = 44251
$ws.Range("K28").Value = 'Angeleno'
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 9500
$ws.Range("S28").Value = 594

# Row 29
$ws.Range("D29").Value = 44251
$ws.Range("K29").Value = 'Angeleno'

# Row 30
$ws.Range("D30").Value = 44279
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 9500
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("S30").Value = 528
$ws.Range("T30").Value = 18

# Row 31
$ws.Range("D31").Value = 44279
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("Q31").Value = '$/bandeja 18 kilos granel'
$ws.Range("S31").Value = 444
$ws.Range("T31").Value = 18

# Row 32
$ws.Range("D32").Value = 44215
$ws.Range("K32").Value = 'Black Amber'
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 11000
$ws.Range("P32").Value = 10500
$ws.Range("S32").Value = 656

# Row 33
$ws.Range("D33").Value = 44215
$ws.Range("K33").Value = 'Black Amber'

# Row 34
$ws.Range("D34").Value = 44257
$ws.Range("N34").Value = 10000
$ws.Range("O34").Value = 11000
$ws.Range("P34").Value = 10500
$ws.Range("Q34").Value = '$/caja 15 kilos granel'
$ws.Range("S34").Value = 700
$ws.Range("T34").Value = 15

# Row 35
$ws.Range("D35").Value = 44257
$ws.Range("N35").Value = 9000
$ws.Range("O35").Value = 9000
$ws.Range("P35").Value = 9000
$ws.Range("Q35").Value = '$/caja 15 kilos granel'
$ws.Range("S35").Value = 600
$ws.Range("T35").Value = 15

# Row 36
$ws.Range("D36").Value = 44323
$ws.Range("K36").Value = 'Angeleno'
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 11000
$ws.Range("O36").Value = 12000
$ws.Range("P36").Value = 11500
$ws.Range("Q36").Value = '$/bandeja 18 kilos granel'
$ws.Range("S36").Value = 639
$ws.Range("T36").Value = 18

# Row 37
$ws.Range("D37").Value = 44323
$ws.Range("K37").Value = 'Angeleno'
$ws.Range("M37").Value = 100
$ws.Range("Q37").Value = '$/bandeja 18 kilos granel'
$ws.Range("S37").Value = 500
$ws.Range("T37").Value = 18
